# Refresh currentAveragePrice / LevePrice / LeveProfit figures pulled by
# the scheduled market-data runner, across all per-job leve sheets.
$wb = $excel.ActiveWorkbook

$rowUpdates = @(
    @{ sheet="ALC"; set=@{ "H76"=8650.695; "I76"=9212.083000000001; "K76"=9212.083000000001; "M76"=-8897.083000000001 } },
    @{ sheet="ALC"; set=@{ "H79"=8650.695; "I79"=9212.083000000001; "K79"=9212.083000000001; "M79"=-8120.083000000001 } },
    @{ sheet="ALC"; set=@{ "H112"=3071.05; "I112"=1373; "J112"=3495.5625; "K112"=4119; "L112"=10486.6875; "M112"=-3011; "N112"=-12702.6875 } },
    @{ sheet="ALC"; set=@{ "H137"=2889.4666; "I137"=2725.5; "K137"=8176.5; "M137"=-5626.5 } },
    @{ sheet="ALC"; set=@{ "H138"=3408.1292; "I138"=1965.1428; "J138"=3829; "K138"=5895.428400000001; "L138"=11487; "M138"=-755.4284000000007; "N138"=-21767 } },
    @{ sheet="ARM"; set=@{ "H2"=3692.1035; "I2"=3664.4614; "K2"=3664.4614; "M2"=-3551.4614 } },
    @{ sheet="ARM"; set=@{ "H32"=8834.045; "I32"=764.94446; "K32"=764.94446; "M32"=-477.94446 } },
    @{ sheet="ARM"; set=@{ "H45"=9888.839; "I45"=13877.444; "J45"=4366.154; "K45"=13877.444; "L45"=4366.154; "M45"=-13500.444; "N45"=-5120.154 } },
    @{ sheet="ARM"; set=@{ "H63"=3436.3; "I63"=3436.3; "K63"=3436.3; "M63"=-2750.3 } },
    @{ sheet="ARM"; set=@{ "H66"=3436.3; "I66"=3436.3; "K66"=17181.5; "M66"=-13749.5 } },
    @{ sheet="ARM"; set=@{ "H116"=3692.1035; "I116"=3664.4614; "K116"=3664.4614; "M116"=-1370.4614 } },
    @{ sheet="BSM"; set=@{ "H3"=3692.1035; "I3"=3664.4614; "K3"=3664.4614; "M3"=-3550.4614 } },
    @{ sheet="BSM"; set=@{ "H20"=1707.381; "I20"=2054; "J20"=1014.1429; "K20"=2054; "L20"=1014.1429; "M20"=-1807; "N20"=-1508.1429 } },
    @{ sheet="CRP"; set=@{ "H8"=1000; "I8"=1000; "J8"=0; "K8"=1000; "L8"=0; "M8"=-860 }; clear=@("N8") },
    @{ sheet="CRP"; set=@{ "H31"=5364.391; "I31"=6632.625; "K31"=6632.625; "M31"=-6337.625 } },
    @{ sheet="CRP"; set=@{ "H34"=5364.391; "I34"=6632.625; "K34"=6632.625; "M34"=-6430.625 } },
    @{ sheet="CRP"; set=@{ "H35"=2198.5; "I35"=2198.5; "K35"=2198.5; "M35"=-1904.5 } },
    @{ sheet="CRP"; set=@{ "H137"=67389.5; "J137"=69999; "L137"=69999; "N137"=-80199 } },
    @{ sheet="CRP"; set=@{ "H141"=86359.60000000001; "J141"=86359.60000000001; "L141"=86359.60000000001; "N141"=-96719.60000000001 } },
    @{ sheet="CUL"; set=@{ "H4"=44198108; "I4"=67380000; "K4"=202140000; "M4"=-202139888 } },
    @{ sheet="CUL"; set=@{ "H18"=325.42856; "I18"=325.42856; "K18"=976.28568; "M18"=-807.28568 } },
    @{ sheet="CUL"; set=@{ "H36"=1500; "I36"=1500; "K36"=4500; "M36"=-4331 } },
    @{ sheet="CUL"; set=@{ "H50"=612.1177; "I50"=348.9; "J50"=988.1429000000001; "K50"=1046.7; "L50"=2964.4287; "M50"=-565.6999999999998; "N50"=-3926.4287 } },
    @{ sheet="CUL"; set=@{ "H53"=612.1177; "I53"=348.9; "J53"=988.1429000000001; "K53"=1046.7; "L53"=2964.4287; "M53"=-565.6999999999998; "N53"=-3926.4287 } },
    @{ sheet="CUL"; set=@{ "H62"=7000; "I62"=6250; "K62"=18750; "M62"=-18064 } },
    @{ sheet="CUL"; set=@{ "H65"=7000; "I65"=6250; "K65"=56250; "M65"=-52818 } },
    @{ sheet="CUL"; set=@{ "H112"=6416.394; "I112"=3457.5; "J112"=7073.926; "K112"=10372.5; "L112"=21221.778; "M112"=-9264.5; "N112"=-23437.778 } },
    @{ sheet="GSM"; set=@{ "H132"=5212.1; "I132"=5360.4287; "J132"=4866; "K132"=16081.2861; "L132"=14598; "M132"=-13551.2861; "N132"=-19658 } },
    @{ sheet="LTW"; set=@{ "H22"=2480.5715; "I22"=2343.5; "K22"=2343.5; "M22"=-2048.5 } },
    @{ sheet="LTW"; set=@{ "H27"=2480.5715; "I27"=2343.5; "K27"=2343.5; "M27"=-2236.5 } },
    @{ sheet="LTW"; set=@{ "H40"=14232.8125; "I40"=11157.625; "J40"=17308; "K40"=11157.625; "L40"=17308; "M40"=-11021.625; "N40"=-17580 } },
    @{ sheet="LTW"; set=@{ "H61"=3791.2222; "I61"=3557.7693; "J61"=4398.2; "K61"=3557.7693; "L61"=4398.2; "M61"=-3355.7693; "N61"=-4802.2 } },
    @{ sheet="LTW"; set=@{ "H93"=1097.5; "I93"=1097.5; "K93"=1097.5; "M93"=150.5 } },
    @{ sheet="LTW"; set=@{ "H113"=3791.2222; "I113"=3557.7693; "J113"=4398.2; "K113"=3557.7693; "L113"=4398.2; "M113"=-1387.7693; "N113"=-8738.200000000001 } },
    @{ sheet="LTW"; set=@{ "H122"=6139.273; "I122"=5042.25; "J122"=6766.143; "K122"=15126.75; "L122"=20298.429; "M122"=-12676.75; "N122"=-25198.429 } },
    @{ sheet="WVR"; set=@{ "H62"=5400; "I62"=5625; "K62"=5625; "M62"=-5001 } },
    @{ sheet="WVR"; set=@{ "H65"=5400; "I65"=5625; "K65"=28125; "M65"=-25005 } }
)

foreach ($u in $rowUpdates) {
    $ws = $wb.Worksheets.Item($u.sheet)
    foreach ($cellRef in $u.set.Keys) {
        $ws.Range($cellRef).Value = $u.set[$cellRef]
    }
    if ($u.ContainsKey("clear")) {
        foreach ($cellRef in $u.clear) {
            $ws.Range($cellRef).ClearContents()
        }
    }
}